$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 13
$ws.Range("I3").Value = 13
$ws.Range("I4").Value = 12
$ws.Range("I5").Value = 13
